$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ------------------------------------------------------------------
# Append the new validation row (row 68) with its raw date/time values
# ------------------------------------------------------------------
$ws.Range("A68").Value = 42474.875
$ws.Range("B68").Value = 42474.875
$ws.Range("C68").Value = 42474.895833333336
$ws.Range("D68").Value = 42474.895833333336
$ws.Range("E68").Value = 42474.895833333336
$ws.Range("F68").Value = 42474.895833333336

# C68:F68 just need the normal date/time display format (reuses the
# existing number format already used throughout the sheet).
$ws.Range("C68:F68").NumberFormat = "dd/mm/yy\ hh:mm"

# B68 reuses the workbook's existing thick red "flag" border (copy the
# formatting from an existing cell that already carries it so the style
# table entry is reused rather than duplicated).
$ws.Range("A6").Copy()
$ws.Range("B68").PasteSpecial(-4122)
$ws.Range("B68").Value = 42474.875

# A68 gets a brand new thick blue border style - the "extra validation
# check" flagged for this row. Start from the existing thick-border
# style (so numbering/weight/format match) and then recolor it blue.
$ws.Range("A6").Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("A68").Value = 42474.875
$ws.Range("A68").Borders.Color = 16711680

$excel.CutCopyMode = 0

$ws.Range("A69").Select
